# Update column F ("dSF") values for specific rows as per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 5
    4  = 3
    6  = 12
    7  = 1
    8  = 4
    9  = 1
    10 = -2
    12 = 3
    14 = -2
    18 = 0
    19 = 10
    20 = -10
    22 = 9
    24 = 3
    25 = -1
    26 = -1
    30 = -5
    31 = 0
    32 = -1
    33 = 1
    34 = 1
    36 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
